# Remove the "reviews_count" column (column E) from the header row,
# shifting all subsequent header values one column to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
